# Regenerate save_data to use K (Strikeouts) instead of Strike# in column G.
# Update column G ("K") values for rows 2-28 on Sheet1 with newly computed values.

$wb = $excel.ActiveWorkbook
$ws = $wb.Sheets.Item("Sheet1")

$kValues = @{
    2  = 0
    3  = 2
    4  = 1
    5  = 1
    6  = 0
    7  = 0
    8  = 4
    9  = 1
    10 = 2
    11 = 0
    12 = 3
    13 = 3
    14 = 2
    15 = 1
    16 = 2
    17 = 2
    18 = 3
    19 = 4
    20 = 4
    21 = 5
    22 = 3
    23 = 6
    24 = 3
    25 = 3
    26 = 5
    27 = 1
    28 = 1
}

foreach ($row in $kValues.Keys) {
    $ws.Range("G$row").Value = $kValues[$row]
}
